$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 3426.4211
$ws.Cells.Item(76, 9).Value = 3306.1875
$ws.Cells.Item(76, 10).Value = 4067.6667
$ws.Cells.Item(76, 11).Value = 3306.1875
$ws.Cells.Item(76, 12).Value = 4067.6667
$ws.Cells.Item(76, 13).Value = -2991.1875
$ws.Cells.Item(76, 14).Value = -4697.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(79, 8).Value = 3426.4211
$ws.Cells.Item(79, 9).Value = 3306.1875
$ws.Cells.Item(79, 10).Value = 4067.6667
$ws.Cells.Item(79, 11).Value = 3306.1875
$ws.Cells.Item(79, 12).Value = 4067.6667
$ws.Cells.Item(79, 13).Value = -2214.1875
$ws.Cells.Item(79, 14).Value = -6251.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 3209.4211
$ws.Cells.Item(98, 9).Value = 2659.8462
$ws.Cells.Item(98, 10).Value = 4400.1665
$ws.Cells.Item(98, 11).Value = 2659.8462
$ws.Cells.Item(98, 12).Value = 4400.1665
$ws.Cells.Item(98, 13).Value = -1161.8462
$ws.Cells.Item(98, 14).Value = -7396.1665

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 3209.4211
$ws.Cells.Item(122, 9).Value = 2659.8462
$ws.Cells.Item(122, 10).Value = 4400.1665
$ws.Cells.Item(122, 11).Value = 7979.5386
$ws.Cells.Item(122, 12).Value = 13200.4995
$ws.Cells.Item(122, 13).Value = -5529.5386
$ws.Cells.Item(122, 14).Value = -18100.4995

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 6457752
$ws.Cells.Item(132, 9).Value = 8006232
$ws.Cells.Item(132, 10).Value = 5752
$ws.Cells.Item(132, 11).Value = 24018696
$ws.Cells.Item(132, 12).Value = 17256
$ws.Cells.Item(132, 13).Value = -24016166
$ws.Cells.Item(132, 14).Value = -22316

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 533.1053000000001
$ws.Cells.Item(135, 9).Value = 568.2941
$ws.Cells.Item(135, 10).Value = 234
$ws.Cells.Item(135, 11).Value = 5114.6469
$ws.Cells.Item(135, 12).Value = 2106
$ws.Cells.Item(135, 13).Value = -2579.6469
$ws.Cells.Item(135, 14).Value = -7176

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 5004643
$ws.Cells.Item(137, 9).Value = 6254804
$ws.Cells.Item(137, 10).Value = 4000.75
$ws.Cells.Item(137, 11).Value = 18764412
$ws.Cells.Item(137, 12).Value = 12002.25
$ws.Cells.Item(137, 13).Value = -18761862
$ws.Cells.Item(137, 14).Value = -17102.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1375.9474
$ws.Cells.Item(61, 9).Value = 988.9677
$ws.Cells.Item(61, 10).Value = 3089.7144
$ws.Cells.Item(61, 11).Value = 988.9677
$ws.Cells.Item(61, 12).Value = 3089.7144
$ws.Cells.Item(61, 13).Value = -776.9677
$ws.Cells.Item(61, 14).Value = -3513.7144

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 1198.7273
$ws.Cells.Item(74, 9).Value = 910.2857
$ws.Cells.Item(74, 10).Value = 1703.5
$ws.Cells.Item(74, 11).Value = 910.2857
$ws.Cells.Item(74, 12).Value = 1703.5
$ws.Cells.Item(74, 13).Value = -36.28570000000002
$ws.Cells.Item(74, 14).Value = -3451.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 1198.7273
$ws.Cells.Item(77, 9).Value = 910.2857
$ws.Cells.Item(77, 10).Value = 1703.5
$ws.Cells.Item(77, 11).Value = 4551.4285
$ws.Cells.Item(77, 12).Value = 8517.5
$ws.Cells.Item(77, 13).Value = -183.4285
$ws.Cells.Item(77, 14).Value = -17253.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 13).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(91, 8).Value = 0
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 11).Value = 0
$ws.Cells.Item(91, 12).Value = 0
$ws.Cells.Item(91, 13).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 25004460
$ws.Cells.Item(132, 9).Value = 40004420
$ws.Cells.Item(132, 10).Value = 4528.4
$ws.Cells.Item(132, 11).Value = 120013260
$ws.Cells.Item(132, 12).Value = 13585.2
$ws.Cells.Item(132, 13).Value = -120010730
$ws.Cells.Item(132, 14).Value = -18645.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 1375.9474
$ws.Cells.Item(136, 9).Value = 988.9677
$ws.Cells.Item(136, 10).Value = 3089.7144
$ws.Cells.Item(136, 11).Value = 2966.9031
$ws.Cells.Item(136, 12).Value = 9269.143199999999
$ws.Cells.Item(136, 13).Value = -416.9031
$ws.Cells.Item(136, 14).Value = -14369.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2093.8
$ws.Cells.Item(105, 9).Value = 1840
$ws.Cells.Item(105, 10).Value = 2202.5715
$ws.Cells.Item(105, 11).Value = 1840
$ws.Cells.Item(105, 12).Value = 2202.5715
$ws.Cells.Item(105, 13).Value = -93
$ws.Cells.Item(105, 14).Value = -5696.5715

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 3726.0386
$ws.Cells.Item(134, 9).Value = 3498.9048
$ws.Cells.Item(134, 10).Value = 4680
$ws.Cells.Item(134, 11).Value = 10496.7144
$ws.Cells.Item(134, 12).Value = 14040
$ws.Cells.Item(134, 13).Value = -7961.714399999999
$ws.Cells.Item(134, 14).Value = -19110

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3034040.5
$ws.Cells.Item(31, 9).Value = 4169055.5
$ws.Cells.Item(31, 10).Value = 7333.3335
$ws.Cells.Item(31, 11).Value = 4169055.5
$ws.Cells.Item(31, 12).Value = 7333.3335
$ws.Cells.Item(31, 13).Value = -4168760.5
$ws.Cells.Item(31, 14).Value = -7923.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 3034040.5
$ws.Cells.Item(34, 9).Value = 4169055.5
$ws.Cells.Item(34, 10).Value = 7333.3335
$ws.Cells.Item(34, 11).Value = 4169055.5
$ws.Cells.Item(34, 12).Value = 7333.3335
$ws.Cells.Item(34, 13).Value = -4168853.5
$ws.Cells.Item(34, 14).Value = -7737.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 16131266
$ws.Cells.Item(58, 9).Value = 924.1111
$ws.Cells.Item(58, 10).Value = 38465584
$ws.Cells.Item(58, 11).Value = 924.1111
$ws.Cells.Item(58, 12).Value = 38465584
$ws.Cells.Item(58, 13).Value = -721.1111
$ws.Cells.Item(58, 14).Value = -38465990

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 3814.3103
$ws.Cells.Item(132, 9).Value = 2665
$ws.Cells.Item(132, 10).Value = 5998
$ws.Cells.Item(132, 11).Value = 7995
$ws.Cells.Item(132, 12).Value = 17994
$ws.Cells.Item(132, 13).Value = -5465
$ws.Cells.Item(132, 14).Value = -23054

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 1560.2593
$ws.Cells.Item(134, 9).Value = 929.2353000000001
$ws.Cells.Item(134, 10).Value = 2633
$ws.Cells.Item(134, 11).Value = 2787.7059
$ws.Cells.Item(134, 12).Value = 7899
$ws.Cells.Item(134, 13).Value = -252.7058999999999
$ws.Cells.Item(134, 14).Value = -12969

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 16131266
$ws.Cells.Item(136, 9).Value = 924.1111
$ws.Cells.Item(136, 10).Value = 38465584
$ws.Cells.Item(136, 11).Value = 2772.3333
$ws.Cells.Item(136, 12).Value = 115396752
$ws.Cells.Item(136, 13).Value = -222.3332999999998
$ws.Cells.Item(136, 14).Value = -115401852

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 44044.906
$ws.Cells.Item(131, 9).Value = 2085
$ws.Cells.Item(131, 10).Value = 53917.824
$ws.Cells.Item(131, 11).Value = 6255
$ws.Cells.Item(131, 12).Value = 161753.472
$ws.Cells.Item(131, 13).Value = -1215
$ws.Cells.Item(131, 14).Value = -171833.472

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(4, 8).Value = 60337
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 60337
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 60337
$ws.Cells.Item(4, 14).Value = -60561

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 7607.3076
$ws.Cells.Item(5, 9).Value = 3000
$ws.Cells.Item(5, 10).Value = 10486.875
$ws.Cells.Item(5, 11).Value = 3000
$ws.Cells.Item(5, 12).Value = 10486.875
$ws.Cells.Item(5, 13).Value = -2888
$ws.Cells.Item(5, 14).Value = -10710.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(100, 8).Value = 30000
$ws.Cells.Item(100, 9).Value = 0
$ws.Cells.Item(100, 10).Value = 30000
$ws.Cells.Item(100, 11).Value = 0
$ws.Cells.Item(100, 12).Value = 30000
$ws.Cells.Item(100, 14).Value = -32164

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2047.9166
$ws.Cells.Item(102, 9).Value = 1458.0358
$ws.Cells.Item(102, 10).Value = 4112.5
$ws.Cells.Item(102, 11).Value = 1458.0358
$ws.Cells.Item(102, 12).Value = 4112.5
$ws.Cells.Item(102, 13).Value = 163.9641999999999
$ws.Cells.Item(102, 14).Value = -7356.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3520.9312
$ws.Cells.Item(132, 9).Value = 3007.2666
$ws.Cells.Item(132, 10).Value = 4071.2856
$ws.Cells.Item(132, 11).Value = 9021.799800000001
$ws.Cells.Item(132, 12).Value = 12213.8568
$ws.Cells.Item(132, 13).Value = -6491.799800000001
$ws.Cells.Item(132, 14).Value = -17273.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 3033.9678
$ws.Cells.Item(132, 9).Value = 1758.1538
$ws.Cells.Item(132, 10).Value = 3955.389
$ws.Cells.Item(132, 11).Value = 5274.4614
$ws.Cells.Item(132, 12).Value = 11866.167
$ws.Cells.Item(132, 13).Value = -2744.4614
$ws.Cells.Item(132, 14).Value = -16926.167

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 2704659.5
$ws.Cells.Item(136, 9).Value = 4001284.2
$ws.Cells.Item(136, 10).Value = 3357.5
$ws.Cells.Item(136, 11).Value = 12003852.6
$ws.Cells.Item(136, 12).Value = 10072.5
$ws.Cells.Item(136, 13).Value = -12001302.6
$ws.Cells.Item(136, 14).Value = -15172.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 283505.56
$ws.Cells.Item(132, 9).Value = 372877.84
$ws.Cells.Item(132, 10).Value = 15388.667
$ws.Cells.Item(132, 11).Value = 1118633.52
$ws.Cells.Item(132, 12).Value = 46166.001
$ws.Cells.Item(132, 13).Value = -1116103.52
$ws.Cells.Item(132, 14).Value = -51226.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 1754.3334
$ws.Cells.Item(136, 9).Value = 865
$ws.Cells.Item(136, 10).Value = 4200
$ws.Cells.Item(136, 11).Value = 2595
$ws.Cells.Item(136, 12).Value = 12600
$ws.Cells.Item(136, 13).Value = -45
$ws.Cells.Item(136, 14).Value = -17700
